$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset (x, y, z) for data rows 2..31 after the edit.
# Rows 2-3 are newly inserted, rows 4-23 are the original rows 2-21 shifted
# down by two, and rows 24-31 are newly appended at the end.
$data = @(
    @(-0.7247905336188185,  -24.30215097461226,   -0.4334598116620501),
    @(0.9856911530861703,   -27.7952582342386,    -10.58349611914382),
    @(-1.217567012860227,   -15.8785987887862,     2.25199632955022),
    @(0.9202353192504198,   -13.56062749151649,   -2.717396190180546),
    @(11.59089219358541,    -11.35514022330562,   -0.3842343388224623),
    @(22.46794862013592,    -12.58473284569012,    5.044853869274501),
    @(3.344646702151348,    -15.18862899802845,   -5.453918620679495),
    @(-2.965254491602879,   -16.62569452319625,    4.355243773150008),
    @(3.611347477816996,    -15.73287558132374,   -0.8439290142623932),
    @(-2.01500372914868,    -25.76902237164192,   -10.39586599902989),
    @(11.75571148212126,    -16.33359923729541,    0.1073646178611298),
    @(39.60853858812318,    -45.04104654887719,    25.44283379887685),
    @(-21.39290519578826,   -25.75901896854879,    6.399868327485425),
    @(-34.85155020945177,   -10.35765630959052,   -4.708384362903586),
    @(-1.052037713090453,   -12.47812661616762,   -0.7732271075954033),
    @(22.74343358553365,    -12.66562087154951,    10.2766472348094),
    @(21.01473594981556,    -12.33190507719506,    14.17957940750583),
    @(0.4384320727466777,   -30.00794573225216,    5.117792562620172),
    @(1.731652502477468,    -23.00902230076561,   -15.75403659865687),
    @(-8.11227344761256,      5.674290752974979,  -28.99636612841348),
    @(-38.96337933512132,   -34.5937611280813,     3.551583487606536),
    @(-22.12481265378448,   -36.7449983704023,     14.93422636336839),
    @(1.603906876942109,     -6.532943116137262,   5.020883938264563),
    @(-0.676244735717717,   -14.42456348125731,    5.205562151395338),
    @(-8.344637724069424,   -29.37389576646717,    4.180637551482669),
    @(-18.19668616345639,    -0.7442300390208842, -9.97852061345025),
    @(-38.98753061238127,   -23.6091353738096,    -9.86452339668957),
    @(-29.63654207760074,   -38.7297830186655,    21.44472347208752),
    @(-12.67720028352449,    -6.331371307373074,   2.991470156336736),
    @(15.31067461092843,    -20.30018748898473,   31.22906570886048)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
